# Fix list level numbering: decrement the outline/indent level of the
# affected paragraphs by one (PowerPoint's IndentLevel is 1-based, so this
# corresponds to lowering the OOXML <a:pPr lvl="..."/> value by one).

$p = $ppt.ActivePresentation

# Slide 1 - "Content Placeholder 2": bulleted list paragraphs.
$s1 = $p.Slides.Item(1)
$tr1 = $s1.Shapes.Item(2).TextFrame.TextRange
$tr1.Paragraphs(1, 1).IndentLevel = 1  # "Bulleted bulleted lists."   lvl 1 -> 0
$tr1.Paragraphs(2, 1).IndentLevel = 1  # "And go to arbitrary depth." lvl 1 -> 0
$tr1.Paragraphs(3, 1).IndentLevel = 2  # "Like this"                  lvl 2 -> 1
$tr1.Paragraphs(4, 1).IndentLevel = 3  # "Or this"                    lvl 3 -> 2
$tr1.Paragraphs(5, 1).IndentLevel = 2  # "Back to here."              lvl 2 -> 1

# Slide 2 - "Content Placeholder 2": numbered list paragraphs.
$s2 = $p.Slides.Item(2)
$tr2 = $s2.Shapes.Item(2).TextFrame.TextRange
$tr2.Paragraphs(2, 1).IndentLevel = 1  # "Tomatoes"                       lvl 1 -> 0
$tr2.Paragraphs(3, 1).IndentLevel = 1  # "Potatoes of various sorts"      lvl 1 -> 0
$tr2.Paragraphs(4, 1).IndentLevel = 2  # "sweet potatoes"                 lvl 2 -> 1
$tr2.Paragraphs(5, 1).IndentLevel = 2  # "russet potates"                 lvl 2 -> 1
$tr2.Paragraphs(6, 1).IndentLevel = 1  # "Tornadoes, for the rhyme."      lvl 1 -> 0
